$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.915632
$ws.Range("H2").Value = 35.746896
$ws.Range("I2").Value = 0.2203762099850903
$ws.Range("J2").Value = 0.2203762099850904
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.126464333333334
$ws.Range("N2").Value = 18.379393
$ws.Range("O2").Value = 0.1081098818071741
$ws.Range("P2").Value = 0.1081098818071741
$ws.Range("Q2").Value = 73.00069445712533
$ws.Range("R2").Value = 657.006250114128
$ws.Range("S2").Value = 0.0238248460146011
$ws.Range("T2").Value = 0.02382484601460111

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.915632
$ws.Range("H3").Value = 35.746896
$ws.Range("I3").Value = 0.2203762099850903
$ws.Range("J3").Value = 0.2203762099850904
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 35.73736866666667
$ws.Range("N3").Value = 107.212106
$ws.Range("O3").Value = 0.6306349784216607
$ws.Range("P3").Value = 0.6306349784216608
$ws.Range("Q3").Value = 425.8333336803307
$ws.Range("R3").Value = 3832.500003122976
$ws.Range("S3").Value = 0.1389769464285948
$ws.Range("T3").Value = 0.1389769464285948

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.915632
$ws.Range("H4").Value = 35.746896
$ws.Range("I4").Value = 0.2203762099850903
$ws.Range("J4").Value = 0.2203762099850904
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.80503233333333
$ws.Range("N4").Value = 44.415097
$ws.Range("O4").Value = 0.2612551397711651
$ws.Range("P4").Value = 0.2612551397711651
$ws.Range("Q4").Value = 176.4113170321014
$ws.Range("R4").Value = 1587.701853288912
$ws.Range("S4").Value = 0.05757441754189441
$ws.Range("T4").Value = 0.05757441754189441

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 29.800487
$ws.Range("H5").Value = 89.40146100000001
$ws.Range("I5").Value = 0.5511514941691683
$ws.Range("J5").Value = 0.5511514941691684
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.126464333333334
$ws.Range("N5").Value = 18.379393
$ws.Range("O5").Value = 0.1081098818071741
$ws.Range("P5").Value = 0.1081098818071741
$ws.Range("Q5").Value = 182.5716207214637
$ws.Range("R5").Value = 1643.144586493173
$ws.Range("S5").Value = 0.0595849228924762
$ws.Range("T5").Value = 0.05958492289247622

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 29.800487
$ws.Range("H6").Value = 89.40146100000001
$ws.Range("I6").Value = 0.5511514941691683
$ws.Range("J6").Value = 0.5511514941691684
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 35.73736866666667
$ws.Range("N6").Value = 107.212106
$ws.Range("O6").Value = 0.6306349784216607
$ws.Range("P6").Value = 0.6306349784216608
$ws.Range("Q6").Value = 1064.990990365208
$ws.Range("R6").Value = 9584.918913286867
$ws.Range("S6").Value = 0.3475754106324395
$ws.Range("T6").Value = 0.3475754106324396

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 29.800487
$ws.Range("H7").Value = 89.40146100000001
$ws.Range("I7").Value = 0.5511514941691683
$ws.Range("J7").Value = 0.5511514941691684
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.80503233333333
$ws.Range("N7").Value = 44.415097
$ws.Range("O7").Value = 0.2612551397711651
$ws.Range("P7").Value = 0.2612551397711651
$ws.Range("Q7").Value = 441.1971735840798
$ws.Range("R7").Value = 3970.774562256718
$ws.Range("S7").Value = 0.1439911606442525
$ws.Range("T7").Value = 0.1439911606442526

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.35338333333333
$ws.Range("H8").Value = 37.06015
$ws.Range("I8").Value = 0.2284722958457413
$ws.Range("J8").Value = 0.2284722958457413
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.126464333333334
$ws.Range("N8").Value = 18.379393
$ws.Range("O8").Value = 0.1081098818071741
$ws.Range("P8").Value = 0.1081098818071741
$ws.Range("Q8").Value = 75.68256238766112
$ws.Range("R8").Value = 681.14306148895
$ws.Range("S8").Value = 0.02470011290009681
$ws.Range("T8").Value = 0.02470011290009681

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.35338333333333
$ws.Range("H9").Value = 37.06015
$ws.Range("I9").Value = 0.2284722958457413
$ws.Range("J9").Value = 0.2284722958457413
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 35.73736866666667
$ws.Range("N9").Value = 107.212106
$ws.Range("O9").Value = 0.6306349784216607
$ws.Range("P9").Value = 0.6306349784216608
$ws.Range("Q9").Value = 441.4774144639889
$ws.Range("R9").Value = 3973.2967301759
$ws.Range("S9").Value = 0.1440826213606263
$ws.Range("T9").Value = 0.1440826213606264

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.35338333333333
$ws.Range("H10").Value = 37.06015
$ws.Range("I10").Value = 0.2284722958457413
$ws.Range("J10").Value = 0.2284722958457413
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.80503233333333
$ws.Range("N10").Value = 44.415097
$ws.Range("O10").Value = 0.2612551397711651
$ws.Range("P10").Value = 0.2612551397711651
$ws.Range("Q10").Value = 182.8922396760611
$ws.Range("R10").Value = 1646.03015708455
$ws.Range("S10").Value = 0.05968956158501812
$ws.Range("T10").Value = 0.05968956158501813

